$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 stays empty (blank row separator), touch it minimally so it
# materializes as an empty <row r="54"/> element without extra attrs.
$ws.Rows.Item(54).OutlineLevel = 0

# Row 55 gets the new report entry.
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "2025-04-28 03:36:03"
$ws.Range("C55").Value = "James Davis shipped New Battery from Ford to Suppliers Old Reliable.`nNow James Davis is Confident, feeling that the task was Challenging.`n"

$ws.Range("C55").WrapText = $true
$ws.Rows.Item(55).AutoFit()
